$wb = $excel.ActiveWorkbook

# The access-control rules live on the first worksheet ("Sheet1"), which
# is already the active sheet in this workbook.
$ws = $wb.Sheets.Item("Sheet1")

# Fill in the two new rule rows (28 = Complaint, 29 = Case File) that
# deny access to non-participants when the "restricted" flag is set.
$ws.Range("B28").Value = "Complaint – Restricted Flag"
$ws.Range("C28").Value = "COMPLAINT"
$ws.Range("D28").Value = "restricted"
$ws.Range("G28").Value = "deny read to *"

$ws.Range("B29").Value = "Case File – Restricted Flag"
$ws.Range("C29").Value = "CASE_FILE"
$ws.Range("D29").Value = "restricted"
$ws.Range("G29").Value = "deny read to *"

# Update the view so the newly added rows are visible, matching the
# author's final cursor/scroll position when the change was made.
$ws.Range("A7").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B30").Select()
